# "Modif ligne en trop horaire masso.xlsx"
# The MASSO sheet had a stray extra row of placeholder ("TBD") schedule
# entries in B9:C9 that don't belong (row 9 should only carry the D column
# note). Clear that leftover line, then leave MASSO as the active/selected
# sheet with that cleared range selected (matching where the author was
# last working).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MASSO")

# Remove the extra leftover "TBD" cells in row 9 (columns B and C).
$ws.Range("B9:C9").Clear()

# MASSO becomes the active sheet/tab, with B9:C9 left selected.
$ws.Activate()
[void]$ws.Range("B9:C9").Select()
